$d = $word.ActiveDocument

# --- Change 1: remove the standalone "Meta description" paragraph that
# currently sits right after the title heading (it reappears, reshaped,
# near the end of the document in Change 2). ---
$metaIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Meta description*") {
        $metaIdx = $i
        break
    }
}
if ($metaIdx -eq -1) {
    throw "Could not locate the 'Meta description' paragraph"
}
[void]$d.Paragraphs.Item($metaIdx).Range.Delete()

# --- Change 2: find the closing "Prompt: ..." image-prompt paragraph
# (now the last paragraph in the document) and, right before it, insert a
# new bold "Play Break da Bank Again Respins Free | Slot Review"
# paragraph; then replace the prompt paragraph's own text with the
# (former) meta-description sentence, keeping its italic formatting. ---
$promptIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Prompt:*") {
        $promptIdx = $i
        break
    }
}
if ($promptIdx -eq -1) {
    throw "Could not locate the 'Prompt: ...' paragraph"
}

$promptPara = $d.Paragraphs.Item($promptIdx)
$rng = $d.Range($promptPara.Range.Start, $promptPara.Range.End)

$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$titleParaXml = "<w:p $w><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Break da Bank Again Respins Free | Slot Review</w:t></w:r></w:p>"
$descParaXml = "<w:p $w><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Discover the respins mechanic and intuitive gameplay of Break da Bank Again Respins by Microgaming. Play free and read our review to learn more.</w:t></w:r></w:p>"

[void]$rng.InsertXML($titleParaXml + $descParaXml)
